# Applies arrear schedule correction ("container related warning rectified")
# The annual increment (Old Basic / New Basic, col C/D) is shifted to align
# correctly so each increment block now starts one row/month later; the
# dependent DA-linked columns (F, G, H) are recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 78300
$ws.Cells.Item(2, 4).Value = 81200
$ws.Cells.Item(2, 6).Value = 86130
$ws.Cells.Item(2, 7).Value = 89320
$ws.Cells.Item(2, 8).Value = 3190

$ws.Cells.Item(3, 3).Value = 78300
$ws.Cells.Item(3, 4).Value = 81200
$ws.Cells.Item(3, 6).Value = 86130
$ws.Cells.Item(3, 7).Value = 89320
$ws.Cells.Item(3, 8).Value = 3190

$ws.Cells.Item(4, 3).Value = 78300
$ws.Cells.Item(4, 4).Value = 81200
$ws.Cells.Item(4, 6).Value = 86130
$ws.Cells.Item(4, 7).Value = 89320
$ws.Cells.Item(4, 8).Value = 3190

$ws.Cells.Item(5, 3).Value = 78300
$ws.Cells.Item(5, 4).Value = 81200
$ws.Cells.Item(5, 6).Value = 86130
$ws.Cells.Item(5, 7).Value = 89320
$ws.Cells.Item(5, 8).Value = 3190

$ws.Cells.Item(6, 3).Value = 78300
$ws.Cells.Item(6, 4).Value = 81200
$ws.Cells.Item(6, 6).Value = 86130
$ws.Cells.Item(6, 7).Value = 89320
$ws.Cells.Item(6, 8).Value = 3190

$ws.Cells.Item(7, 3).Value = 78300
$ws.Cells.Item(7, 4).Value = 81200
$ws.Cells.Item(7, 6).Value = 86130
$ws.Cells.Item(7, 7).Value = 89320
$ws.Cells.Item(7, 8).Value = 3190

$ws.Cells.Item(8, 3).Value = 78300
$ws.Cells.Item(8, 4).Value = 81200
$ws.Cells.Item(8, 6).Value = 86130
$ws.Cells.Item(8, 7).Value = 89320
$ws.Cells.Item(8, 8).Value = 3190

$ws.Cells.Item(9, 3).Value = 78300
$ws.Cells.Item(9, 4).Value = 81200
$ws.Cells.Item(9, 6).Value = 86130
$ws.Cells.Item(9, 7).Value = 89320
$ws.Cells.Item(9, 8).Value = 3190

$ws.Cells.Item(10, 3).Value = 78300
$ws.Cells.Item(10, 4).Value = 81200
$ws.Cells.Item(10, 6).Value = 86130
$ws.Cells.Item(10, 7).Value = 89320
$ws.Cells.Item(10, 8).Value = 3190

$ws.Cells.Item(11, 3).Value = 78300
$ws.Cells.Item(11, 4).Value = 81200
$ws.Cells.Item(11, 6).Value = 86130
$ws.Cells.Item(11, 7).Value = 89320
$ws.Cells.Item(11, 8).Value = 3190

$ws.Cells.Item(12, 3).Value = 78300
$ws.Cells.Item(12, 4).Value = 81200
$ws.Cells.Item(12, 6).Value = 86130
$ws.Cells.Item(12, 7).Value = 89320
$ws.Cells.Item(12, 8).Value = 3190

$ws.Cells.Item(13, 3).Value = 78300
$ws.Cells.Item(13, 4).Value = 81200
$ws.Cells.Item(13, 6).Value = 86130
$ws.Cells.Item(13, 7).Value = 89320
$ws.Cells.Item(13, 8).Value = 3190

$ws.Cells.Item(14, 3).Value = 80700
$ws.Cells.Item(14, 4).Value = 83700
$ws.Cells.Item(14, 6).Value = 91191
$ws.Cells.Item(14, 7).Value = 94581
$ws.Cells.Item(14, 8).Value = 3390

$ws.Cells.Item(15, 3).Value = 80700
$ws.Cells.Item(15, 4).Value = 83700
$ws.Cells.Item(15, 6).Value = 91191
$ws.Cells.Item(15, 7).Value = 94581
$ws.Cells.Item(15, 8).Value = 3390

$ws.Cells.Item(16, 3).Value = 80700
$ws.Cells.Item(16, 4).Value = 83700
$ws.Cells.Item(16, 6).Value = 91191
$ws.Cells.Item(16, 7).Value = 94581
$ws.Cells.Item(16, 8).Value = 3390

$ws.Cells.Item(17, 3).Value = 80700
$ws.Cells.Item(17, 4).Value = 83700
$ws.Cells.Item(17, 6).Value = 91191
$ws.Cells.Item(17, 7).Value = 94581
$ws.Cells.Item(17, 8).Value = 3390

$ws.Cells.Item(18, 3).Value = 80700
$ws.Cells.Item(18, 4).Value = 83700
$ws.Cells.Item(18, 6).Value = 91191
$ws.Cells.Item(18, 7).Value = 94581
$ws.Cells.Item(18, 8).Value = 3390

$ws.Cells.Item(19, 3).Value = 80700
$ws.Cells.Item(19, 4).Value = 83700
$ws.Cells.Item(19, 6).Value = 91191
$ws.Cells.Item(19, 7).Value = 94581
$ws.Cells.Item(19, 8).Value = 3390

$ws.Cells.Item(20, 3).Value = 80700
$ws.Cells.Item(20, 4).Value = 83700
$ws.Cells.Item(20, 6).Value = 91191
$ws.Cells.Item(20, 7).Value = 94581
$ws.Cells.Item(20, 8).Value = 3390

$ws.Cells.Item(21, 3).Value = 80700
$ws.Cells.Item(21, 4).Value = 83700
$ws.Cells.Item(21, 6).Value = 91191
$ws.Cells.Item(21, 7).Value = 94581
$ws.Cells.Item(21, 8).Value = 3390

$ws.Cells.Item(22, 3).Value = 80700
$ws.Cells.Item(22, 4).Value = 83700
$ws.Cells.Item(22, 6).Value = 91191
$ws.Cells.Item(22, 7).Value = 94581
$ws.Cells.Item(22, 8).Value = 3390

$ws.Cells.Item(23, 3).Value = 80700
$ws.Cells.Item(23, 4).Value = 83700
$ws.Cells.Item(23, 6).Value = 91191
$ws.Cells.Item(23, 7).Value = 94581
$ws.Cells.Item(23, 8).Value = 3390

$ws.Cells.Item(24, 3).Value = 80700
$ws.Cells.Item(24, 4).Value = 83700
$ws.Cells.Item(24, 6).Value = 91191
$ws.Cells.Item(24, 7).Value = 94581
$ws.Cells.Item(24, 8).Value = 3390

$ws.Cells.Item(25, 3).Value = 80700
$ws.Cells.Item(25, 4).Value = 83700
$ws.Cells.Item(25, 6).Value = 91191
$ws.Cells.Item(25, 7).Value = 94581
$ws.Cells.Item(25, 8).Value = 3390

$ws.Cells.Item(26, 3).Value = 83200
$ws.Cells.Item(26, 4).Value = 86300
$ws.Cells.Item(26, 6).Value = 94016
$ws.Cells.Item(26, 7).Value = 97519
$ws.Cells.Item(26, 8).Value = 3503

$ws.Cells.Item(27, 3).Value = 83200
$ws.Cells.Item(27, 4).Value = 86300
$ws.Cells.Item(27, 6).Value = 94016
$ws.Cells.Item(27, 7).Value = 97519
$ws.Cells.Item(27, 8).Value = 3503

$ws.Cells.Item(28, 3).Value = 83200
$ws.Cells.Item(28, 4).Value = 86300
$ws.Cells.Item(28, 6).Value = 94016
$ws.Cells.Item(28, 7).Value = 97519
$ws.Cells.Item(28, 8).Value = 3503

$ws.Cells.Item(29, 3).Value = 83200
$ws.Cells.Item(29, 4).Value = 86300
$ws.Cells.Item(29, 6).Value = 94016
$ws.Cells.Item(29, 7).Value = 97519
$ws.Cells.Item(29, 8).Value = 3503

$ws.Cells.Item(30, 3).Value = 83200
$ws.Cells.Item(30, 4).Value = 86300
$ws.Cells.Item(30, 6).Value = 94016
$ws.Cells.Item(30, 7).Value = 97519
$ws.Cells.Item(30, 8).Value = 3503

$ws.Cells.Item(31, 3).Value = 83200
$ws.Cells.Item(31, 4).Value = 86300
$ws.Cells.Item(31, 6).Value = 94016
$ws.Cells.Item(31, 7).Value = 97519
$ws.Cells.Item(31, 8).Value = 3503

$ws.Cells.Item(32, 3).Value = 83200
$ws.Cells.Item(32, 4).Value = 86300
$ws.Cells.Item(32, 6).Value = 94016
$ws.Cells.Item(32, 7).Value = 97519
$ws.Cells.Item(32, 8).Value = 3503

$ws.Cells.Item(33, 3).Value = 83200
$ws.Cells.Item(33, 4).Value = 86300
$ws.Cells.Item(33, 6).Value = 94016
$ws.Cells.Item(33, 7).Value = 97519
$ws.Cells.Item(33, 8).Value = 3503

$ws.Cells.Item(34, 3).Value = 83200
$ws.Cells.Item(34, 4).Value = 86300
$ws.Cells.Item(34, 6).Value = 94016
$ws.Cells.Item(34, 7).Value = 97519
$ws.Cells.Item(34, 8).Value = 3503

$ws.Cells.Item(35, 3).Value = 83200
$ws.Cells.Item(35, 4).Value = 86300
$ws.Cells.Item(35, 6).Value = 94016
$ws.Cells.Item(35, 7).Value = 97519
$ws.Cells.Item(35, 8).Value = 3503

$ws.Cells.Item(36, 3).Value = 83200
$ws.Cells.Item(36, 4).Value = 86300
$ws.Cells.Item(36, 6).Value = 94016
$ws.Cells.Item(36, 7).Value = 97519
$ws.Cells.Item(36, 8).Value = 3503

$ws.Cells.Item(37, 3).Value = 83200
$ws.Cells.Item(37, 4).Value = 86300
$ws.Cells.Item(37, 6).Value = 94016
$ws.Cells.Item(37, 7).Value = 97519
$ws.Cells.Item(37, 8).Value = 3503

$ws.Cells.Item(38, 3).Value = 85700
$ws.Cells.Item(38, 4).Value = 88900
$ws.Cells.Item(38, 6).Value = 96841
$ws.Cells.Item(38, 7).Value = 100457
$ws.Cells.Item(38, 8).Value = 3616

$ws.Cells.Item(39, 3).Value = 85700
$ws.Cells.Item(39, 4).Value = 88900
$ws.Cells.Item(39, 6).Value = 96841
$ws.Cells.Item(39, 7).Value = 100457
$ws.Cells.Item(39, 8).Value = 3616

$ws.Cells.Item(40, 3).Value = 85700
$ws.Cells.Item(40, 4).Value = 88900
$ws.Cells.Item(40, 6).Value = 99412
$ws.Cells.Item(40, 7).Value = 103124
$ws.Cells.Item(40, 8).Value = 3712

$ws.Cells.Item(41, 3).Value = 85700
$ws.Cells.Item(41, 4).Value = 88900
$ws.Cells.Item(41, 6).Value = 99412
$ws.Cells.Item(41, 7).Value = 103124
$ws.Cells.Item(41, 8).Value = 3712

$ws.Cells.Item(42, 3).Value = 85700
$ws.Cells.Item(42, 4).Value = 88900
$ws.Cells.Item(42, 6).Value = 99412
$ws.Cells.Item(42, 7).Value = 103124
$ws.Cells.Item(42, 8).Value = 3712

$ws.Cells.Item(43, 3).Value = 85700
$ws.Cells.Item(43, 4).Value = 88900
$ws.Cells.Item(43, 6).Value = 99412
$ws.Cells.Item(43, 7).Value = 103124
$ws.Cells.Item(43, 8).Value = 3712

$ws.Cells.Item(44, 3).Value = 85700
$ws.Cells.Item(44, 4).Value = 88900
$ws.Cells.Item(44, 6).Value = 99412
$ws.Cells.Item(44, 7).Value = 103124
$ws.Cells.Item(44, 8).Value = 3712

$ws.Cells.Item(45, 3).Value = 85700
$ws.Cells.Item(45, 4).Value = 88900
$ws.Cells.Item(45, 6).Value = 99412
$ws.Cells.Item(45, 7).Value = 103124
$ws.Cells.Item(45, 8).Value = 3712

$ws.Cells.Item(46, 2).Value = 6600
$ws.Cells.Item(46, 3).Value = 85700
$ws.Cells.Item(46, 4).Value = 88900
$ws.Cells.Item(46, 6).Value = 99412
$ws.Cells.Item(46, 7).Value = 103124
$ws.Cells.Item(46, 8).Value = 3712

$ws.Cells.Item(47, 2).Value = 6600
$ws.Cells.Item(47, 3).Value = 85700
$ws.Cells.Item(47, 4).Value = 88900
$ws.Cells.Item(47, 6).Value = 99412
$ws.Cells.Item(47, 7).Value = 103124
$ws.Cells.Item(47, 8).Value = 3712

$ws.Cells.Item(48, 2).Value = 6600
$ws.Cells.Item(48, 3).Value = 85700
$ws.Cells.Item(48, 4).Value = 88900
$ws.Cells.Item(48, 6).Value = 99412
$ws.Cells.Item(48, 7).Value = 103124
$ws.Cells.Item(48, 8).Value = 3712

$ws.Cells.Item(49, 2).Value = 6600
$ws.Cells.Item(49, 3).Value = 85700
$ws.Cells.Item(49, 4).Value = 88900
$ws.Cells.Item(49, 6).Value = 99412
$ws.Cells.Item(49, 7).Value = 103124
$ws.Cells.Item(49, 8).Value = 3712

$ws.Cells.Item(58, 3).Value = 99800
$ws.Cells.Item(58, 4).Value = 105700
$ws.Cells.Item(58, 6).Value = 123752
$ws.Cells.Item(58, 7).Value = 131068
$ws.Cells.Item(58, 8).Value = 7316

$ws.Cells.Item(59, 3).Value = 99800
$ws.Cells.Item(59, 4).Value = 105700
$ws.Cells.Item(59, 6).Value = 123752
$ws.Cells.Item(59, 7).Value = 131068
$ws.Cells.Item(59, 8).Value = 7316

$ws.Cells.Item(60, 3).Value = 99800
$ws.Cells.Item(60, 4).Value = 105700
$ws.Cells.Item(60, 6).Value = 123752
$ws.Cells.Item(60, 7).Value = 131068
$ws.Cells.Item(60, 8).Value = 7316

$ws.Cells.Item(61, 3).Value = 99800
$ws.Cells.Item(61, 4).Value = 105700
$ws.Cells.Item(61, 6).Value = 123752
$ws.Cells.Item(61, 7).Value = 131068
$ws.Cells.Item(61, 8).Value = 7316

$ws.Cells.Item(70, 3).Value = 102800
$ws.Cells.Item(70, 4).Value = 108900
$ws.Cells.Item(70, 6).Value = 131584
$ws.Cells.Item(70, 7).Value = 139392
$ws.Cells.Item(70, 8).Value = 7808

$ws.Cells.Item(71, 3).Value = 102800
$ws.Cells.Item(71, 4).Value = 108900
$ws.Cells.Item(71, 6).Value = 131584
$ws.Cells.Item(71, 7).Value = 139392
$ws.Cells.Item(71, 8).Value = 7808

$ws.Cells.Item(72, 3).Value = 102800
$ws.Cells.Item(72, 4).Value = 108900
$ws.Cells.Item(72, 6).Value = 131584
$ws.Cells.Item(72, 7).Value = 139392
$ws.Cells.Item(72, 8).Value = 7808

$ws.Cells.Item(73, 3).Value = 102800
$ws.Cells.Item(73, 4).Value = 108900
$ws.Cells.Item(73, 6).Value = 131584
$ws.Cells.Item(73, 7).Value = 139392
$ws.Cells.Item(73, 8).Value = 7808

